$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 20:05"
$ws.Range("B4").Value = 2577515
$ws.Range("C4").Value = 24559
$ws.Range("D4").Value = 1071393
$ws.Range("E4").Value = 1378170
$ws.Range("G4").Value = 312
$ws.Range("H4").Value = 127952
$ws.Range("B5").Value = 1284214
$ws.Range("C5").Value = 4160
$ws.Range("E5").Value = 530491
$ws.Range("G5").Value = 88
$ws.Range("H5").Value = 56197
$ws.Range("B7").Value = 529331
$ws.Range("C7").Value = 19885
$ws.Range("D7").Value = 310120
$ws.Range("E7").Value = 203109
$ws.Range("G7").Value = 413
$ws.Range("H7").Value = 16102
$ws.Range("B9").Value = 295549
$ws.Range("C9").Value = 564
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 28341
$ws.Range("A16").Value = "Turquia"
$ws.Range("B16").Value = 195883
$ws.Range("C16").Value = 1372
$ws.Range("D16").Value = 169182
$ws.Range("E16").Value = 21619
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 5082
$ws.Range("A17").Value = "Alemania"
$ws.Range("B17").Value = 194539
$ws.Range("C17").Value = 140
$ws.Range("D17").Value = 177500
$ws.Range("E17").Value = 8013
$ws.Range("H17").Value = 9026
$ws.Range("B34").Value = 47360
$ws.Range("C34").Value = 387
$ws.Range("D34").Value = 35834
$ws.Range("E34").Value = 11215
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 311
$ws.Range("A57").Value = "Ghana"
$ws.Range("B57").Value = 16431
$ws.Range("C57").Value = 597
$ws.Range("D57").Value = 12257
$ws.Range("E57").Value = 4071
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 103
$ws.Range("A58").Value = "Moldavia"
$ws.Range("B58").Value = 16080
$ws.Range("C58").Value = 304
$ws.Range("D58").Value = 8963
$ws.Range("E58").Value = 6596
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 521
$ws.Range("A59").Value = "Honduras"
$ws.Range("B59").Value = 15994
$ws.Range("C59").Value = 628
$ws.Range("D59").Value = 1678
$ws.Range("E59").Value = 13845
$ws.Range("G59").Value = 45
$ws.Range("H59").Value = 471
$ws.Range("A60").Value = "Azerbaiyan"
$ws.Range("B60").Value = 15890
$ws.Range("C60").Value = 521
$ws.Range("D60").Value = 8719
$ws.Range("E60").Value = 6978
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 193
$ws.Range("B68").Value = 11877
$ws.Range("C68").Value = 244
$ws.Range("D68").Value = 8723
$ws.Range("E68").Value = 2934
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 220
$ws.Range("B75").Value = 7551
$ws.Range("C75").Value = 124
$ws.Range("D75").Value = 5240
$ws.Range("E75").Value = 2291
$ws.Range("A95").Value = "Guayana Francesa"
$ws.Range("B95").Value = 3461
$ws.Range("C95").Value = 191
$ws.Range("D95").Value = 1249
$ws.Range("E95").Value = 2200
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 12
$ws.Range("A96").Value = "Grecia"
$ws.Range("B96").Value = 3366
$ws.Range("C96").Value = 23
$ws.Range("D96").Value = 1374
$ws.Range("E96").Value = 1801
$ws.Range("H96").Value = 191
$ws.Range("A97").Value = "Republica de Africa Central"
$ws.Range("B97").Value = 3340
$ws.Range("D97").Value = 661
$ws.Range("E97").Value = 2639
$ws.Range("H97").Value = 40
$ws.Range("E100").Value = 1543
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 13
$ws.Range("B105").Value = 2305
$ws.Range("C105").Value = 22
$ws.Range("D105").Value = 1875
$ws.Range("E105").Value = 422
$ws.Range("B107").Value = 2118
$ws.Range("C107").Value = 58
$ws.Range("D107").Value = 1398
$ws.Range("E107").Value = 607
$ws.Range("A114").Value = "Estado de Palestina"
$ws.Range("B114").Value = 1815
$ws.Range("C114").Value = 258
$ws.Range("D114").Value = 446
$ws.Range("E114").Value = 1366
$ws.Range("H114").Value = 3
$ws.Range("A115").Value = "Lituania"
$ws.Range("B115").Value = 1813
$ws.Range("C115").Value = 5
$ws.Range("D115").Value = 1503
$ws.Range("E115").Value = 232
$ws.Range("H115").Value = 78
$ws.Range("A116").Value = "Libano"
$ws.Range("B116").Value = 1719
$ws.Range("C116").Value = 22
$ws.Range("D116").Value = 1153
$ws.Range("E116").Value = 533
$ws.Range("H116").Value = 33
$ws.Range("A117").Value = "Paraguay"
$ws.Range("B117").Value = 1711
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 1013
$ws.Range("E117").Value = 685
$ws.Range("H117").Value = 13
$ws.Range("A118").Value = "Eslovaquia"
$ws.Range("B118").Value = 1657
$ws.Range("C118").Value = 14
$ws.Range("D118").Value = 1455
$ws.Range("E118").Value = 174
$ws.Range("H118").Value = 28
$ws.Range("B128").Value = 1103
$ws.Range("C128").Value = 14
$ws.Range("D128").Value = 417
$ws.Range("E128").Value = 390
$ws.Range("G128").Value = 3
$ws.Range("H128").Value = 296
$ws.Range("A155").Value = "Montenegro"
$ws.Range("B155").Value = 469
$ws.Range("C155").Value = 30
$ws.Range("D155").Value = 315
$ws.Range("E155").Value = 145
$ws.Range("H155").Value = 9
$ws.Range("A156").Value = "Taiwan"
$ws.Range("B156").Value = 447
$ws.Range("D156").Value = 435
$ws.Range("E156").Value = 5
$ws.Range("H156").Value = 7
$ws.Range("A164").Value = "Angola"
$ws.Range("B164").Value = 244
$ws.Range("C164").Value = 32
$ws.Range("D164").Value = 81
$ws.Range("E164").Value = 153
$ws.Range("H164").Value = 10
$ws.Range("A165").Value = "Martinica"
$ws.Range("B165").Value = 242
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 98
$ws.Range("E165").Value = 130
$ws.Range("H165").Value = 14
$ws.Range("A166").Value = "Guyana"
$ws.Range("B166").Value = 230
$ws.Range("C166").Value = 15
$ws.Range("D166").Value = 109
$ws.Range("E166").Value = 109
$ws.Range("H166").Value = 12
$ws.Range("A167").Value = "Mongolia"
$ws.Range("B167").Value = 219
$ws.Range("D167").Value = 175
$ws.Range("E167").Value = 44
$ws.Range("H167").Value = 0
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Laos"
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A204").Value = "Dominica"
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"
